# Add four new Saudi-city rows (149-152) to Sheet1, mirroring the
# append-only edit described by the upstream diff (new rows only;
# the first 148 rows are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, in the same column layout as the existing table:
#   A: City (Arabic-alphabet City field, actually holds the English name
#      in this sheet, mirroring column B)
#   B: City English Name
#   C: City Arabic Name
#   D: Latitude
#   E: Longitude
#   F: Region (Arabic)
#   G: Area (Arabic)
$newRows = @(
    @("Al Namas",  "Al Namas",  "النماص", 19.115787999999998, 42.168045999999997, "منطقة عسير", "جنوب المملكة"),
    @("Billasmar", "Billasmar", "بللسمر", 18.793106000000002, 42.255910999999998, "منطقة عسير", "جنوب المملكة"),
    @("Al Harth",  "Al Harth",  "الحرث",  16.811212999999999, 43.152147999999997, "منطقة جازان", "جنوب المملكة"),
    @("Baqaa",     "Baqaa",     "بقعاء",  27.889050000000001, 42.415892999999997, "منطقة حائل", "شمال المملكة")
)

$startRow = 149
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy formatting (the thin-border style used by every data row) from
    # the last existing row down onto the new one before filling values.
    $ws.Range("A148:G148").Copy()
    $ws.Range("A$r`:G$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$lastRow = $startRow + $newRows.Count - 1
$ws.Application.CutCopyMode = $false

# Refresh the used-range selection/dimension to span the newly added rows.
$ws.Range("A1:G$lastRow").Select()
